$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new blank columns before column B (old B:M shifts right to E:P,
# all formulas referencing those columns shift automatically).
$ws.Columns("B:D").Insert()

# Rename the existing "Snoods" header (now shifted to K1) to "Snoods/Bibs".
$ws.Range("K1").Value = "Snoods/Bibs"

# New bottom "legend" rows under the table (order chosen to match the
# resulting shared-string table ordering).
$ws.Range("E15").Value = "Missing"

# New header for inserted column D.
$ws.Range("D1").Value = "Diapers"

$ws.Range("E16").Value = "Belly Bands"

# New headers for inserted columns B and C.
$ws.Range("B1").Value = "Bibs"
$ws.Range("C1").Value = "Bellybands"

# Center (but don't bold) the three new headers.
$ws.Range("B1:D1").HorizontalAlignment = -4108
$ws.Range("B1:D1").ColumnWidth = 18.8333333333333

$ws.Range("E17").Value = "Bibs"
$ws.Range("E18").Value = "Diapers"

# The new "Hats" column (shifted from old K to new N) previously had no
# value for the xxs/xs/s rows; the edit fills those in with "NA" to match
# the other "Quote on Request"-style rows.
$ws.Range("N3").Value = "NA"
$ws.Range("N4").Value = "NA"
$ws.Range("N5").Value = "NA"

# Move the active selection to reflect where the user ended up editing.
$ws.Range("E19").Select()
